# Update NATMI LR-pair TPM output values (Ccl11-Ccr5) per new TPM data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.57077
$ws.Range("H2").Value = 4.71231
$ws.Range("I2").Value = 0.02582502173444737
$ws.Range("J2").Value = 0.02582502173444737
$ws.Range("M2").Value = 0.243056
$ws.Range("N2").Value = 0.729168
$ws.Range("O2").Value = 0.002199620488481675
$ws.Range("P2").Value = 0.002199620488481675
$ws.Range("Q2").Value = 0.3817850731200001
$ws.Range("R2").Value = 3.43606565808
$ws.Range("S2").Value = 0.00005680524692257499
$ws.Range("T2").Value = 0.00005680524692257499
$ws.Range("G3").Value = 1.57077
$ws.Range("H3").Value = 4.71231
$ws.Range("I3").Value = 0.02582502173444737
$ws.Range("J3").Value = 0.02582502173444737
$ws.Range("M3").Value = 70.95253000000001
$ws.Range("N3").Value = 212.85759
$ws.Range("O3").Value = 0.6421097964979703
$ws.Range("P3").Value = 0.6421097964979703
$ws.Range("Q3").Value = 111.4501055481
$ws.Range("R3").Value = 1003.0509499329
$ws.Range("S3").Value = 0.01658249945046166
$ws.Range("T3").Value = 0.01658249945046166
$ws.Range("G4").Value = 1.57077
$ws.Range("H4").Value = 4.71231
$ws.Range("I4").Value = 0.02582502173444737
$ws.Range("J4").Value = 0.02582502173444737
$ws.Range("M4").Value = 0.04794200000000001
$ws.Range("N4").Value = 0.143826
$ws.Range("O4").Value = 0.0004338679376719292
$ws.Range("P4").Value = 0.0004338679376719292
$ws.Range("Q4").Value = 0.07530585534000002
$ws.Range("R4").Value = 0.6777526980600002
$ws.Range("S4").Value = 0.00001120464892025743
$ws.Range("T4").Value = 0.00001120464892025743
$ws.Range("G5").Value = 1.57077
$ws.Range("H5").Value = 4.71231
$ws.Range("I5").Value = 0.02582502173444737
$ws.Range("J5").Value = 0.02582502173444737
$ws.Range("M5").Value = 39.25553366666666
$ws.Range("N5").Value = 117.766601
$ws.Range("O5").Value = 0.3552567150758761
$ws.Range("P5").Value = 0.3552567150758761
$ws.Range("Q5").Value = 61.66141461759
$ws.Range("R5").Value = 554.95273155831
$ws.Range("S5").Value = 0.009174512388142876
$ws.Range("T5").Value = 0.009174512388142878
$ws.Range("I6").Value = 0.934831682683009
$ws.Range("J6").Value = 0.934831682683009
$ws.Range("M6").Value = 0.243056
$ws.Range("N6").Value = 0.729168
$ws.Range("O6").Value = 0.002199620488481675
$ws.Range("P6").Value = 0.002199620488481675
$ws.Range("Q6").Value = 13.82011546778134
$ws.Range("R6").Value = 124.381039210032
$ws.Range("S6").Value = 0.002056274922511346
$ws.Range("T6").Value = 0.002056274922511346
$ws.Range("I7").Value = 0.934831682683009
$ws.Range("J7").Value = 0.934831682683009
$ws.Range("M7").Value = 70.95253000000001
$ws.Range("N7").Value = 212.85759
$ws.Range("O7").Value = 0.6421097964979703
$ws.Range("P7").Value = 0.6421097964979703
$ws.Range("Q7").Value = 4034.346641643158
$ws.Range("R7").Value = 36309.11977478841
$ws.Range("S7").Value = 0.600264581527442
$ws.Range("T7").Value = 0.600264581527442
$ws.Range("I8").Value = 0.934831682683009
$ws.Range("J8").Value = 0.934831682683009
$ws.Range("M8").Value = 0.04794200000000001
$ws.Range("N8").Value = 0.143826
$ws.Range("O8").Value = 0.0004338679376719292
$ws.Range("P8").Value = 0.0004338679376719292
$ws.Range("Q8").Value = 2.725972515619334
$ws.Range("R8").Value = 24.533752640574
$ws.Range("S8").Value = 0.0004055934942360565
$ws.Range("T8").Value = 0.0004055934942360565
$ws.Range("I9").Value = 0.934831682683009
$ws.Range("J9").Value = 0.934831682683009
$ws.Range("M9").Value = 39.25553366666666
$ws.Range("N9").Value = 117.766601
$ws.Range("O9").Value = 0.3552567150758761
$ws.Range("P9").Value = 0.3552567150758761
$ws.Range("Q9").Value = 2232.061780094756
$ws.Range("R9").Value = 20088.5560208528
$ws.Range("S9").Value = 0.3321052327388195
$ws.Range("T9").Value = 0.3321052327388196
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.9273763333333335
$ws.Range("H10").Value = 2.782129
$ws.Range("I10").Value = 0.01524698967025436
$ws.Range("J10").Value = 0.01524698967025436
$ws.Range("M10").Value = 0.243056
$ws.Range("N10").Value = 0.729168
$ws.Range("O10").Value = 0.002199620488481675
$ws.Range("P10").Value = 0.002199620488481675
$ws.Range("Q10").Value = 0.2254043820746667
$ws.Range("R10").Value = 2.028639438672
$ws.Range("S10").Value = 0.00003353759086635994
$ws.Range("T10").Value = 0.00003353759086635994
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.9273763333333335
$ws.Range("H11").Value = 2.782129
$ws.Range("I11").Value = 0.01524698967025436
$ws.Range("J11").Value = 0.01524698967025436
$ws.Range("M11").Value = 70.95253000000001
$ws.Range("N11").Value = 212.85759
$ws.Range("O11").Value = 0.6421097964979703
$ws.Range("P11").Value = 0.6421097964979703
$ws.Range("Q11").Value = 65.79969711212335
$ws.Range("R11").Value = 592.1972740091101
$ws.Range("S11").Value = 0.009790241434373683
$ws.Range("T11").Value = 0.009790241434373683
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.9273763333333335
$ws.Range("H12").Value = 2.782129
$ws.Range("I12").Value = 0.01524698967025436
$ws.Range("J12").Value = 0.01524698967025436
$ws.Range("M12").Value = 0.04794200000000001
$ws.Range("N12").Value = 0.143826
$ws.Range("O12").Value = 0.0004338679376719292
$ws.Range("P12").Value = 0.0004338679376719292
$ws.Range("Q12").Value = 0.04446027617266668
$ws.Range("R12").Value = 0.4001424855540001
$ws.Range("S12").Value = 0.000006615179963938469
$ws.Range("T12").Value = 0.000006615179963938469
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.9273763333333335
$ws.Range("H13").Value = 2.782129
$ws.Range("I13").Value = 0.01524698967025436
$ws.Range("J13").Value = 0.01524698967025436
$ws.Range("M13").Value = 39.25553366666666
$ws.Range("N13").Value = 117.766601
$ws.Range("O13").Value = 0.3552567150758761
$ws.Range("P13").Value = 0.3552567150758761
$ws.Range("Q13").Value = 36.40465287483656
$ws.Range("R13").Value = 327.641875873529
$ws.Range("S13").Value = 0.00541659546505038
$ws.Range("T13").Value = 0.005416595465050381
$ws.Range("G14").Value = 0.7810079999999999
$ws.Range("H14").Value = 2.343024
$ws.Range("I14").Value = 0.0128405486320577
$ws.Range("J14").Value = 0.0128405486320577
$ws.Range("M14").Value = 0.243056
$ws.Range("N14").Value = 0.729168
$ws.Range("O14").Value = 0.002199620488481675
$ws.Range("P14").Value = 0.002199620488481675
$ws.Range("Q14").Value = 0.189828680448
$ws.Range("R14").Value = 1.708458124032
$ws.Range("S14").Value = 0.00002824433385441945
$ws.Range("T14").Value = 0.00002824433385441945
$ws.Range("G15").Value = 0.7810079999999999
$ws.Range("H15").Value = 2.343024
$ws.Range("I15").Value = 0.0128405486320577
$ws.Range("J15").Value = 0.0128405486320577
$ws.Range("M15").Value = 70.95253000000001
$ws.Range("N15").Value = 212.85759
$ws.Range("O15").Value = 0.6421097964979703
$ws.Range("P15").Value = 0.6421097964979703
$ws.Range("Q15").Value = 55.41449355024
$ws.Range("R15").Value = 498.73044195216
$ws.Range("S15").Value = 0.008245042069052858
$ws.Range("T15").Value = 0.008245042069052858
$ws.Range("G16").Value = 0.7810079999999999
$ws.Range("H16").Value = 2.343024
$ws.Range("I16").Value = 0.0128405486320577
$ws.Range("J16").Value = 0.0128405486320577
$ws.Range("M16").Value = 0.04794200000000001
$ws.Range("N16").Value = 0.143826
$ws.Range("O16").Value = 0.0004338679376719292
$ws.Range("P16").Value = 0.0004338679376719292
$ws.Range("Q16").Value = 0.037443085536
$ws.Range("R16").Value = 0.336987769824
$ws.Range("S16").Value = 0.000005571102353566985
$ws.Range("T16").Value = 0.000005571102353566985
$ws.Range("G17").Value = 0.7810079999999999
$ws.Range("H17").Value = 2.343024
$ws.Range("I17").Value = 0.0128405486320577
$ws.Range("J17").Value = 0.0128405486320577
$ws.Range("M17").Value = 39.25553366666666
$ws.Range("N17").Value = 117.766601
$ws.Range("O17").Value = 0.3552567150758761
$ws.Range("P17").Value = 0.3552567150758761
$ws.Range("Q17").Value = 30.65888583793599
$ws.Range("R17").Value = 275.929972541424
$ws.Range("S17").Value = 0.004561691126796851
$ws.Range("T17").Value = 0.004561691126796852
$ws.Range("G18").Value = 0.6846153333333334
$ws.Range("H18").Value = 2.053846
$ws.Range("I18").Value = 0.01125575728023152
$ws.Range("J18").Value = 0.01125575728023152
$ws.Range("M18").Value = 0.243056
$ws.Range("N18").Value = 0.729168
$ws.Range("O18").Value = 0.002199620488481675
$ws.Range("P18").Value = 0.002199620488481675
$ws.Range("Q18").Value = 0.1663998644586667
$ws.Range("R18").Value = 1.497598780128
$ws.Range("S18").Value = 0.00002475839432697402
$ws.Range("T18").Value = 0.00002475839432697402
$ws.Range("G19").Value = 0.6846153333333334
$ws.Range("H19").Value = 2.053846
$ws.Range("I19").Value = 0.01125575728023152
$ws.Range("J19").Value = 0.01125575728023152
$ws.Range("M19").Value = 70.95253000000001
$ws.Range("N19").Value = 212.85759
$ws.Range("O19").Value = 0.6421097964979703
$ws.Range("P19").Value = 0.6421097964979703
$ws.Range("Q19").Value = 48.57518997679334
$ws.Range("R19").Value = 437.1767097911401
$ws.Range("S19").Value = 0.007227432016640008
$ws.Range("T19").Value = 0.007227432016640008
$ws.Range("G20").Value = 0.6846153333333334
$ws.Range("H20").Value = 2.053846
$ws.Range("I20").Value = 0.01125575728023152
$ws.Range("J20").Value = 0.01125575728023152
$ws.Range("M20").Value = 0.04794200000000001
$ws.Range("N20").Value = 0.143826
$ws.Range("O20").Value = 0.0004338679376719292
$ws.Range("P20").Value = 0.0004338679376719292
$ws.Range("Q20").Value = 0.03282182831066667
$ws.Range("R20").Value = 0.295396454796
$ws.Range("S20").Value = 0.000004883512198109852
$ws.Range("T20").Value = 0.000004883512198109852
$ws.Range("G21").Value = 0.6846153333333334
$ws.Range("H21").Value = 2.053846
$ws.Range("I21").Value = 0.01125575728023152
$ws.Range("J21").Value = 0.01125575728023152
$ws.Range("M21").Value = 39.25553366666666
$ws.Range("N21").Value = 117.766601
$ws.Range("O21").Value = 0.3552567150758761
$ws.Range("P21").Value = 0.3552567150758761
$ws.Range("Q21").Value = 26.87494026638289
$ws.Range("R21").Value = 241.874462397446
$ws.Range("S21").Value = 0.003998683357066426
$ws.Range("T21").Value = 0.003998683357066427
